$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [string]$MatchText,
        [string]$BodyXml
    )
    $target = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text -eq $MatchText) {
            $target = $cand
            break
        }
    }
    if ($target -eq $null) {
        Write-Host "WARNING: paragraph not found for match [$MatchText]"
        return
    }
    $pr = $target.Range
    $rng = $d.Range($pr.Start, $pr.End - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $BodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# --- 1) Lecturer: Sdfasdf -> William Martin ---
Set-ParagraphXml "Lecturer: `tSdfasdf`r" '<w:r><w:t xml:space="preserve">Lecturer: </w:t><w:tab/><w:t>William Martin</w:t></w:r>'

# --- 2) Room: DASD -> S1 ---
Set-ParagraphXml "Room:`tDASD`r" '<w:r><w:t>Room:</w:t><w:tab/><w:t>S1</w:t></w:r>'

# --- 3) Phone paragraph: rebuild with multiple runs / formatting (incl. highlighted "Delete if adjunct") ---
$phoneBody = ""
$phoneBody += '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Phone:</w:t></w:r>'
$phoneBody += '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/><w:t>+815-1717+ e</w:t></w:r>'
$phoneBody += '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">xt.#   </w:t></w:r>'
$phoneBody += '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>++</w:t></w:r>'
$phoneBody += '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>Delete if adjunct</w:t></w:r>'
$phoneBody += '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>++</w:t></w:r>'
Set-ParagraphXml "Phone:`t11-22 ext. 33`r" $phoneBody

# --- 4) Email: sd@dfg.hfd -> willy@hotmail.com ---
Set-ParagraphXml "Email: `tsd@dfg.hfd `r" '<w:r><w:t xml:space="preserve">Email: </w:t><w:tab/><w:t xml:space="preserve">willy@hotmail.com </w:t></w:r>'

# --- 5) Contact hour: 12:21 PM -> 12:45 PM ---
Set-ParagraphXml "Contact hour: `t12:21 PM`r" '<w:r><w:t xml:space="preserve">Contact hour: </w:t><w:tab/><w:t>12:45 PM</w:t></w:r>'

# --- 6) Course Coordinator paragraph: collapse 4 formatted runs into one plain run ---
Set-ParagraphXml "Course Coordinator: `t+Lecturer Name+`r" '<w:r><w:t xml:space="preserve">Course Coordinator: </w:t><w:tab/><w:t>William Martin</w:t></w:r>'

Write-Host "Edits applied."
